$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "시프트업" (ShiftUp) right after row 11 (엑셀세라퓨틱스),
# which pushes 이노그리드 and all following rows down by one.
$ws.Rows(12).Insert()

# Fill in the data for the newly inserted row.
$ws.Cells.Item(12, 1).Value = "시프트업"
$ws.Cells.Item(12, 2).Value = "2024.06.03~06.13"
$ws.Cells.Item(12, 3).Value = "47,000~60,000"
$ws.Cells.Item(12, 4).Value = "-"
$ws.Cells.Item(12, 5).Value = "340750"
$ws.Cells.Item(12, 6).Value = "한국투자증권,NH투자증권,신한투자증권"

# The table keeps a fixed number of rows, so the previous last row
# (KB스팩28호, now pushed to row 22) drops off the bottom.
$ws.Rows(22).Delete()
